$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = 'id: urn:oid:1.2.208.176.2.22'
$meta.Range("B4").Value = '1.1.0'
$meta.Range("B9").Value = '2023-07-10T23:08:03+02:00'
$meta.Range("B11").Value = 'No display for ContactDetail'

# --- Concepts sheet updates (definition text simplifications) ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C2").Value = 'Social indsats [Servicelov]'
$concepts.Range("D2").Value = 'Indsats, der består af en eller flere ydelser, som leveres til en eller flere borgere af et eller flere tilbud, og som skal forebygge eller tilgodese behov, som følger af en nedsat funktionsevne, og som gives med henblik på at fremme den enkeltes mulighed for at klare sig selv eller at lette den daglige tilværelse eller forbedre livskvaliteten.'
$concepts.Range("D3").Value = 'Tjeneste, genstand eller beløb, der gives eller modtages'
$concepts.Range("D4").Value = 'Ydelse, der har til formål at afdække en persons situation'
$concepts.Range("D5").Value = 'Udredning, hvor en leverandør afdækker en persons funktionsevne i relation til daglige aktiviteter'
$concepts.Range("D6").Value = 'Ydelse, der har til formål at huse en person.'
$concepts.Range("D7").Value = 'Ophold, der er tilkendt for en ikke-midlertidig periode.'
$concepts.Range("D8").Value = 'Ophold, der er tilkendt midlertidigt.'
$concepts.Range("D10").Value = 'Ydelse, der gennem udviklende eller vedligeholdende aktiviteter har til formål at fremme socialt samvær eller tilbyde et miljøskift.'
$concepts.Range("D11").Value = 'Aktivitet og samvær, hvor en person får stimuleret sine sanser med henblik på evnen til at kunne sortere, organisere og bearbejde sanseindtryk.'
$concepts.Range("D12").Value = 'Aktivitet og samvær, hvor en person har individuelle oplevelser eller oplevelser sammen med andre personer.'
$concepts.Range("D13").Value = 'Aktivitet og samvær, der gennem fysiske bevægelser giver en person mulighed for at bruge sin krop.'
$concepts.Range("D14").Value = 'Aktivitet og samvær, som ved at få en person ind i et fællesskab med andre forebygger eller bryder med social isolation og vedligeholder de sociale relationer.'
$concepts.Range("D15").Value = 'Aktivitet og samvær, hvor en person kan udfolde sig kreativt.'
$concepts.Range("D16").Value = 'Aktivitet og samvær, hvor en person eller gruppe af personer undervises i bestemte færdigheder eller kompetencer.'
$concepts.Range("D17").Value = 'Ydelse, der gennem arbejdslignende, lønnede aktiviteter har til formål at afdække, oparbejde, udvikle eller bevare en persons arbejdsevne og beskæftigelsesrelevante kompetencer.'
$concepts.Range("D18").Value = 'Beskyttet beskæftigelse, der er centreret omkring udførelse af servicerende, understøttende arbejdsfunktioner.'
$concepts.Range("D19").Value = 'Beskyttet beskæftigelse, der har til formål, at en person kan afprøve ressourcer og kompetencer i forhold til et specifikt arbejdsområde eller en given arbejdsplads i en afgrænset periode.'
$concepts.Range("D20").Value = 'Beskyttet beskæftigelse, der er centreret omkring udførelse af enkle, afgrænsede arbejdsopgaver som led i en samlet produktion af konkrete produkter.'
$concepts.Range("D21").Value = 'Ydelse, der gennem motivation, vejledning og støtte har til formål at udvikle eller fastholde en persons funktionsevne og muligheder for selvstændighed og selvbestemmelse i forhold til personens situation.'
$concepts.Range("D22").Value = 'Socialpædagogisk støtte, der retter sig mod aktiviteter, der er en forudsætning for at kunne deltage i samfundslivet.'
$concepts.Range("D23").Value = 'Støtte til samfundsdeltagelse, der retter sig mod en persons mulighed for selvstændigt at færdes uden for hjemmet.'
$concepts.Range("D24").Value = 'Støtte til samfundsdeltagelse, der retter sig mod en persons muligheder for at varetage en uddannelse.'
$concepts.Range("D25").Value = 'Støtte til samfundsdeltagelse, der retter sig mod en persons muligheder for at varetage et arbejde.'
$concepts.Range("D26").Value = 'Støtte til samfundsdeltagelse, der retter sig mod gennemførsel af besøg hos og kontakt til offentlige og private instanser med et specifikt formål.'
$concepts.Range("D27").Value = 'Socialpædagogisk støtte, der retter sig mod et eller flere aspekter af sundhed og indeholder et element af sundhedsfremme.'
$concepts.Range("D28").Value = 'Støtte til sundhed, der retter sig mod livsstilsbetingede forhold.'
$concepts.Range("D29").Value = 'Støtte til sundhed, der retter sig mod hygiejnemæssige opgaver i relation til personen selv.'
$concepts.Range("D30").Value = 'Støtte til sundhed, der retter sig mod en persons psykiske velbefindende.'
$concepts.Range("D31").Value = 'Støtte til sundhed, der retter sig mod den måde, som en persons seksualitet kommer til udtryk på.'
$concepts.Range("D32").Value = 'Støtte til sundhed, der retter sig mod en persons behandling.'
$concepts.Range("D33").Value = 'Socialpædagogisk støtte, der retter sig mod relationer og fællesskaber, hvor samspillet med andre er et centralt element.'
$concepts.Range("D34").Value = 'Støtte til relationer og fællesskaber, der retter sig mod at udvikle, indgå i og bevare relationer til andre mennesker.'
$concepts.Range("D35").Value = 'Støtte til relationer og fællesskaber, der retter sig mod strukturering og håndtering af opgaver, som er en almindelig del af forældreansvaret, og hvor samværet mellem barn og forælder er centralt.'
$concepts.Range("D36").Value = 'Socialpædagogisk støtte, der retter sig imod almindelige praktiske opgaver.'
$concepts.Range("D37").Value = 'Støtte til praktiske opgaver, der har hverdagskarakter og retter sig mod en persons hjem.'
$concepts.Range("D38").Value = 'Støtte til praktiske opgaver, der retter sig mod anskaffelse af og etablering i bolig.'
$concepts.Range("D39").Value = 'Støtte til praktiske opgaver af administrativ eller økonomisk karakter.'

# --- Concepts sheet: append mapset info to "FFB tilstande" definition (row 129) ---
$concepts.Range("D129").Value = "FFB tilstande`nMapset: 58000024148 (FFB tilstande mappet til SNOMED)`nMapset: 338000019145 (OmrHierarkiFFB)`nMapset: 68000024145 (TilstandsrelationerFFB)"
